# Updates the crypto price/volume table to reflect the latest scrape.
# Rows 11 and 12 also swap Polygon/TRON (ranking order changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.971.06'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.892.66'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7367'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.78'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3092'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.41'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06903'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7705'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07949'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.894.89'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.218'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.48'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.970.53'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.10'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.792'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.45'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -5.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007767'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.141.70'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.911'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.298'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.98'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.79'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1272'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -9.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.348'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.535'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.295'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.059'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05094'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.276'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7340'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.716'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01922'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.782'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.289'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.11'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4439'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.932'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.48%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8345'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.632'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.01'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.780'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.045.59'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.42'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '931.98'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -5.44%  '
